# Add the "10000 point" data set + answers to Sheet3, per commit:
#   "Large data set working / 10000 points"
#
# New content (rows 20-30, columns A:B):
#   A20/B20 header labels "10000 data set" / "answers"
#   A21:A30 -> input ids, B21:B30 -> computed answers in scientific format

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Header row
$ws.Range("A20").Value = "10000 data set"
$ws.Range("B20").Value = "answers"

# id / answer pairs (answers use "0.00000000E+00" formatting like the rest of col B)
$rows = @(
    @(6736, "2.1802E+19"),
    @(2771, "2.21423E+19"),
    @(9393, "2.24469E+19"),
    @(4278, "2.28697E+19"),
    @(6494, "2.30881E+19"),
    @(9560, "2.34592E+19"),
    @(8839, "2.36409E+19"),
    @(2262, "2.43383E+19"),
    @(5411, "2.43517E+19"),
    @(4373, "2.43994E+19")
)

$r = 21
foreach ($pair in $rows) {
    $ws.Cells.Item($r, 1).Value = $pair[0]

    $answer = [double]$pair[1]
    $ws.Cells.Item($r, 2).Value = $answer
    $ws.Cells.Item($r, 2).NumberFormat = "0.00000000E+00"

    $r = $r + 1
}

# Matches the updated selection anchor left behind in the saved file
$ws.Range("C27").Select()
